$d = $word.ActiveDocument

# Locate the run that needs to be split/corrected: it runs from the start of
# "the industry data section..." through the trailing period after
# "Technlogies" (the typo).
$startRng = $d.Content
$startRng.Find.Execute("the industry data section of the Fighting Mongoose") | Out-Null
$start = $startRng.Start

$endRng = $d.Content
$endRng.Find.Execute("Technlogies.") | Out-Null
$end = $endRng.End

$target = $d.Range($start, $end)

# Replace that single run with five runs (fixing the "Technlogies" ->
# "Technologies" typo along the way) bracketing "in depth" with
# proofErr gramStart/gramEnd markers, matching Word's own re-flow of the
# paragraph after a manual correction.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r><w:t xml:space="preserve">the industry data section of the Fighting Mongoose' + [char]0x2019 + 's homepage provides more </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>in depth</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> analysis using data from Burning Glass </w:t></w:r>' +
  '<w:r><w:t>Technologies</w:t></w:r>' +
  '<w:r><w:t>.</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
